$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'99.315.19"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "'3.288.47"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'254.41"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").Value = "'622.92"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'1.45"
$ws.Range("E7").Value = "  +21.94%  "
$ws.Range("D8").Value = "'0.400"
$ws.Range("E8").Value = "  +2.72%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "'0.982"
$ws.Range("E10").Value = "  +23.36%  "
$ws.Range("D11").Value = "'3.285.39"
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "'39.74"
$ws.Range("E13").Value = "  +10.24%  "
$ws.Range("D14").Value = "'98.946.03"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "'0.0000248"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "'3.907.55"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").Value = "'5.49"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "'3.289.09"
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("D19").Value = "'3.46"
$ws.Range("E19").Value = "  -4.36%  "
$ws.Range("D20").Value = "'15.26"
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").Value = "'6.34"
$ws.Range("E21").Value = "  +7.84%  "
$ws.Range("D22").Value = "'487.21"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'9.30"
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("D24").Value = "'0.0000202"
$ws.Range("E24").Value = "  -3.06%  "
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").Value = "'0.334"
$ws.Range("E26").Value = "  +37.06%  "
$ws.Range("D27").Value = "'89.08"
$ws.Range("D28").Value = "'12.05"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").Value = "'3.461.11"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  +12.68%  "
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").Value = "'10.42"
$ws.Range("E33").Value = "  +12.69%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").Value = "'27.99"
$ws.Range("E35").Value = "  +2.70%  "
$ws.Range("D36").Value = "'0.479"
$ws.Range("E36").Value = "  +7.01%  "
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("D40").Value = "'24.77"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D41").Value = "'487.61"
$ws.Range("E41").Value = "  -5.33%  "
$ws.Range("D42").Value = "'3.75"
$ws.Range("E42").Value = "  +2.83%  "
$ws.Range("E43").Value = "  -2.92%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -1.59%  "
$ws.Range("E46").Value = "  -5.13%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'158.18"
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'1.94"
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("D49").Value = "'0.850"
$ws.Range("E49").Value = "  +7.08%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "'4.75"
$ws.Range("E50").Value = "  +5.36%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'7.29"
$ws.Range("E51").Value = "  +15.04%  "
